# Scheduled runner update: refresh cached Universalis market-price snapshot
# columns for the affected Leve rows (currentAveragePrice* / LevePrice* / LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: "Morning Glass of Ether" (Ether)
$ws.Range("H15").Value = 146.91
$ws.Range("I15").Value = 146.91
$ws.Range("K15").Value = 440.73
$ws.Range("M15").Value = -271.73

# Row 64: "Forged from the Void" (Void Glue)
$ws.Range("H64").Value = 3532.3076
$ws.Range("I64").Value = 3544
$ws.Range("J64").Value = 3506
$ws.Range("K64").Value = 3544
$ws.Range("L64").Value = 3506
$ws.Range("M64").Value = -3296
$ws.Range("N64").Value = -4002

# Row 67: "Dodging the Draft (L)" (Void Glue)
$ws.Range("H67").Value = 3532.3076
$ws.Range("I67").Value = 3544
$ws.Range("J67").Value = 3506
$ws.Range("K67").Value = 3544
$ws.Range("L67").Value = 3506
$ws.Range("M67").Value = -2686
$ws.Range("N67").Value = -5222

# Row 76: "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 3300
$ws.Range("I76").Value = 2950
$ws.Range("J76").Value = 3650
$ws.Range("K76").Value = 2950
$ws.Range("L76").Value = 3650
$ws.Range("M76").Value = -2635
$ws.Range("N76").Value = -4280

# Row 79: "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 3300
$ws.Range("I79").Value = 2950
$ws.Range("J79").Value = 3650
$ws.Range("K79").Value = 2950
$ws.Range("L79").Value = 3650
$ws.Range("M79").Value = -1858
$ws.Range("N79").Value = -5834

# Row 112: "Making Ends Meet" (Superior Spiritbond Potion)
$ws.Range("H112").Value = 1747.1538
$ws.Range("J112").Value = 1820.5278
$ws.Range("L112").Value = 5461.5834
$ws.Range("N112").Value = -7677.5834

# Row 138: "All-night Crafting" (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 6414369
$ws.Range("I138").Value = 13335933
$ws.Range("J138").Value = 5512.926
$ws.Range("K138").Value = 40007799
$ws.Range("L138").Value = 16538.778
$ws.Range("M138").Value = -40002659
$ws.Range("N138").Value = -26818.778

$ws = $wb.Worksheets.Item("ARM")
# Row 32: "Ingot We Trust" (Steel Ingot)
$ws.Range("H32").Value = 21187.094
$ws.Range("I32").Value = 20057.176
$ws.Range("K32").Value = 20057.176
$ws.Range("M32").Value = -19770.176

# Row 74: "As the Bolt Flies" (Titanium Nugget)
$ws.Range("H74").Value = 4558.8276
$ws.Range("I74").Value = 1204.1818
$ws.Range("J74").Value = 15102
$ws.Range("K74").Value = 1204.1818
$ws.Range("L74").Value = 15102
$ws.Range("M74").Value = -330.1818000000001
$ws.Range("N74").Value = -16850

# Row 77: "Heavy Metal Banned (L)" (Titanium Nugget)
$ws.Range("H77").Value = 4558.8276
$ws.Range("I77").Value = 1204.1818
$ws.Range("J77").Value = 15102
$ws.Range("K77").Value = 6020.909000000001
$ws.Range("L77").Value = 75510
$ws.Range("M77").Value = -1652.909000000001
$ws.Range("N77").Value = -84246

$ws = $wb.Worksheets.Item("BSM")
# Row 94: "High Steal" (High Steel Nugget)
$ws.Range("H94").Value = 1042.7407
$ws.Range("I94").Value = 834.9524
$ws.Range("K94").Value = 834.9524
$ws.Range("M94").Value = -383.9524

# Row 105: "Ingot to Wing It" (Molybdenum Ingot)
$ws.Range("H105").Value = 4598.5713
$ws.Range("I105").Value = 3147.5
$ws.Range("J105").Value = 4840.4165
$ws.Range("K105").Value = 3147.5
$ws.Range("L105").Value = 4840.4165
$ws.Range("M105").Value = -1400.5
$ws.Range("N105").Value = -8334.416499999999

# Row 118: "Cooking for the Future" (Titanbronze Culinary Knife)
$ws.Range("H118").Value = 7900
$ws.Range("J118").Value = 7900
$ws.Range("L118").Value = 7900
$ws.Range("N118").Value = -11214

$ws = $wb.Worksheets.Item("CRP")
# Row 31: "Wall Not Found" (Walnut Lumber)
$ws.Range("H31").Value = 5843.185
$ws.Range("I31").Value = 1933.4546
$ws.Range("J31").Value = 8531.125
$ws.Range("K31").Value = 1933.4546
$ws.Range("L31").Value = 8531.125
$ws.Range("M31").Value = -1638.4546
$ws.Range("N31").Value = -9121.125

# Row 34: "Armoires of the Rich and Famous" (Walnut Lumber)
$ws.Range("H34").Value = 5843.185
$ws.Range("I34").Value = 1933.4546
$ws.Range("J34").Value = 8531.125
$ws.Range("K34").Value = 1933.4546
$ws.Range("L34").Value = 8531.125
$ws.Range("M34").Value = -1731.4546
$ws.Range("N34").Value = -8935.125

# Row 62: "Splinter in the Sewers" (Cedar Lumber)
$ws.Range("H62").Value = 2915.4614
$ws.Range("I62").Value = 2443.5715
$ws.Range("J62").Value = 3466
$ws.Range("K62").Value = 2443.5715
$ws.Range("L62").Value = 3466
$ws.Range("M62").Value = -1819.5715
$ws.Range("N62").Value = -4714

# Row 65: "The Lumber of Their Discontent (L)" (Cedar Lumber)
$ws.Range("H65").Value = 2915.4614
$ws.Range("I65").Value = 2443.5715
$ws.Range("J65").Value = 3466
$ws.Range("K65").Value = 12217.8575
$ws.Range("L65").Value = 17330
$ws.Range("M65").Value = -9097.8575
$ws.Range("N65").Value = -23570

# Row 141: "No Greater Treasure" (Claro Walnut Necklace of Gathering)
$ws.Range("H141").Value = 39304.266
$ws.Range("I141").Value = 38000
$ws.Range("J141").Value = 39397.43
$ws.Range("K141").Value = 38000
$ws.Range("L141").Value = 39397.43
$ws.Range("M141").Value = -32820
$ws.Range("N141").Value = -49757.43

$ws = $wb.Worksheets.Item("CUL")
# Row 136: "Simple Is Hardest" (Spaghetti al Olio e Peperoncino)
$ws.Range("H136").Value = 6285
$ws.Range("J136").Value = 6526.316
$ws.Range("L136").Value = 19578.948
$ws.Range("N136").Value = -29778.948

$ws = $wb.Worksheets.Item("GSM")
# Row 70: "Sky Is the Limit" (Mythrite Ingot)
$ws.Range("H70").Value = 11998
$ws.Range("I70").Value = 13998.1
$ws.Range("J70").Value = 3997.6
$ws.Range("K70").Value = 13998.1
$ws.Range("L70").Value = 3997.6
$ws.Range("M70").Value = -13728.1
$ws.Range("N70").Value = -4537.6

# Row 73: "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws.Range("H73").Value = 11998
$ws.Range("I73").Value = 13998.1
$ws.Range("J73").Value = 3997.6
$ws.Range("K73").Value = 13998.1
$ws.Range("L73").Value = 3997.6
$ws.Range("M73").Value = -13062.1
$ws.Range("N73").Value = -5869.6

# Row 80: "Needs More Prayerbell" (Hardsilver Ingot)
$ws.Range("H80").Value = 2928.8572
$ws.Range("I80").Value = 2720.5
$ws.Range("J80").Value = 3449.75
$ws.Range("K80").Value = 2720.5
$ws.Range("L80").Value = 3449.75
$ws.Range("M80").Value = -1722.5
$ws.Range("N80").Value = -5445.75

# Row 83: "With a Noise That Reaches Heaven (L)" (Hardsilver Ingot)
$ws.Range("H83").Value = 2928.8572
$ws.Range("I83").Value = 2720.5
$ws.Range("J83").Value = 3449.75
$ws.Range("K83").Value = 13602.5
$ws.Range("L83").Value = 17248.75
$ws.Range("M83").Value = -8610.5
$ws.Range("N83").Value = -27232.75

# Row 126: "Gold Rush Order" (Phrygian Gold Ingot)
$ws.Range("H126").Value = 3498.3076
$ws.Range("I126").Value = 2279.6667
$ws.Range("J126").Value = 4542.857
$ws.Range("K126").Value = 6839.000100000001
$ws.Range("L126").Value = 13628.571
$ws.Range("M126").Value = -4369.000100000001
$ws.Range("N126").Value = -18568.571

$ws = $wb.Worksheets.Item("LTW")
# Row 40: "Best Served Toad" (Toad Leather)
$ws.Range("H40").Value = 3886.1365
$ws.Range("I40").Value = 4626.4
$ws.Range("K40").Value = 4626.4
$ws.Range("M40").Value = -4490.4

# Row 136: "Respect for Br'aax" (Br'aax Leather)
$ws.Range("H136").Value = 4444.5747
$ws.Range("I136").Value = 2099.5715
$ws.Range("J136").Value = 11284.167
$ws.Range("K136").Value = 6298.7145
$ws.Range("L136").Value = 33852.501
$ws.Range("M136").Value = -3748.7145
$ws.Range("N136").Value = -38952.501

$ws = $wb.Worksheets.Item("WVR")
# Row 136: "Weaving the Envelope" (Sarcenet Cloth)
$ws.Range("H136").Value = 1332.5873
$ws.Range("I136").Value = 648.2982
$ws.Range("J136").Value = 7833.3335
$ws.Range("K136").Value = 1944.8946
$ws.Range("L136").Value = 23500.0005
$ws.Range("M136").Value = 605.1054000000001
$ws.Range("N136").Value = -28600.0005
